# Update header labels on the VisioData sheet so the "From"/"To" line
# attribute headers and the ShapeLabel FontSize header read with spaces
# between the compound words (e.g. "From LineLabel" -> "From Line Label").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VisioData")

$ws.Range("F1").Value = "Shape Label Font Size"
$ws.Range("X1").Value = "From Line Label"
$ws.Range("Y1").Value = "From Line Pattern"
$ws.Range("Z1").Value = "From Arrow Type"
$ws.Range("AA1").Value = "From Line Color"
$ws.Range("AC1").Value = "To Line Label"
$ws.Range("AD1").Value = "To Line Pattern"
$ws.Range("AE1").Value = "To Arrow Type"
$ws.Range("AF1").Value = "To Line Color"

# Move the active selection on the VisioData sheet to AB1 (matches the
# saved cursor position recorded in the workbook).
$ws.Activate()
$ws.Range("AB1").Select()
